$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.408218
$ws.Range("H2").Value = 61.224654
$ws.Range("I2").Value = 0.1108535210972707
$ws.Range("J2").Value = 0.1108535210972707
$ws.Range("M2").Value = 0.464324
$ws.Range("N2").Value = 1.392972
$ws.Range("Q2").Value = 9.476025414632
$ws.Range("R2").Value = 85.28422873168799
$ws.Range("S2").Value = 0.1108535210972707
$ws.Range("T2").Value = 0.1108535210972707

$ws.Range("I3").Value = 0.2566851044076959
$ws.Range("J3").Value = 0.256685104407696
$ws.Range("M3").Value = 0.464324
$ws.Range("N3").Value = 1.392972
$ws.Range("Q3").Value = 21.942059655376
$ws.Range("R3").Value = 197.478536898384
$ws.Range("S3").Value = 0.2566851044076959
$ws.Range("T3").Value = 0.256685104407696

$ws.Range("G4").Value = 85.307233
$ws.Range("H4").Value = 255.921699
$ws.Range("I4").Value = 0.4633725077375833
$ws.Range("J4").Value = 0.4633725077375833
$ws.Range("M4").Value = 0.464324
$ws.Range("N4").Value = 1.392972
$ws.Range("Q4").Value = 39.61019565549199
$ws.Range("R4").Value = 356.491760899428
$ws.Range("S4").Value = 0.4633725077375833
$ws.Range("T4").Value = 0.4633725077375833

$ws.Range("G5").Value = 31.12938966666666
$ws.Range("H5").Value = 93.38816899999999
$ws.Range("I5").Value = 0.16908886675745
$ws.Range("J5").Value = 0.16908886675745
$ws.Range("M5").Value = 0.464324
$ws.Range("N5").Value = 1.392972
$ws.Range("Q5").Value = 14.45412272758533
$ws.Range("R5").Value = 130.087104548268
$ws.Range("S5").Value = 0.16908886675745
$ws.Range("T5").Value = 0.16908886675745
